# Update the division-problem table with a new set of generated problems.
# The document contains a single 5-column table; only rows 1, 5, 9, 13 and 17
# (1-based) hold the visible "a÷b=" expressions, one per cell, five per row.
# We address each cell directly by (row, column) rather than doing a global
# Find/Replace, because several of the new values coincide with *other*
# cells' old values (e.g. "16÷4=" and "77÷7=" are both a source value in one
# cell and a target value in another), which would make a sequential
# Find/Replace ambiguous/unsafe.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("77÷6=", "94÷6=", "44÷6=", "19÷6=", "33÷6=")
    5  = @("54÷3=", "53÷3=", "77÷9=", "91÷9=", "16÷4=")
    9  = @("44÷6=", "35÷8=", "19÷2=", "43÷9=", "86÷8=")
    13 = @("93÷5=", "52÷8=", "77÷7=", "22÷2=", "55÷5=")
    17 = @("50÷3=", "71÷8=", "62÷9=", "27÷4=", "76÷8=")
}

foreach ($rowIndex in $newValues.Keys) {
    $values = $newValues[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
